$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-record row is inserted at row 242 (a weekly "Betarraga" reading
# for 2022-02-09, i.e. serial date 44609, volume 100). Inserting the row
# shifts the existing rows 242:350 down to 243:351 (so the former row 350
# becomes the new row 351), matching the new sheet dimension A1:R351.
$ws.Rows.Item(242).Insert()

# After the shift, the data that used to live in row 242 is now in row 243.
# Duplicate that row into the freshly inserted (blank) row 242 so every
# column keeps its original value/format, then overwrite just the two
# columns (Fecha/D and Volumen/J) that hold the new record's own data.
$ws.Rows.Item(243).Copy()
$ws.Rows.Item(242).PasteSpecial()

$ws.Range("D242").Value = 44609
$ws.Range("J242").Value = 100
